$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Footer "last-updated" date fields: 9/12/2023 -> 9/15/2023
#    These live on the Slide Master, every Custom (slide) Layout, and
#    the Notes Master - not on the slides themselves.
# ---------------------------------------------------------------------

function Set-DatePlaceholderText($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date*") {
            $shp.TextFrame.TextRange.Text = $newText
        }
    }
}

$master = $p.SlideMaster

# Slide master's own Date placeholder
Set-DatePlaceholderText $master.Shapes "9/15/2023"

# Every slide layout hanging off the master
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Set-DatePlaceholderText $layout.Shapes "9/15/2023"
}

# Notes master date field (updates via the headers/footers object)
$notesMaster = $p.NotesMaster
$notesMaster.HeadersFooters.DateAndTime.Text = "9/15/2023"

# ---------------------------------------------------------------------
# 2) Slide 5, "TextBox 15": refresh the Off-grid / building progress
#    bullets and grow the textbox to fit the new content.
# ---------------------------------------------------------------------

$slide5 = $p.Slides.Item(5)
$shp = $slide5.Shapes.Item("TextBox 15")
$tf = $shp.TextFrame
$tr = $tf.TextRange

function Replace-BulletText($oldText, $newText) {
    $full = $tr.Text
    $idx = $full.IndexOf($oldText)
    if ($idx -ge 0) {
        $sub = $tr.Characters($idx + 1, $oldText.Length)
        $sub.Text = $newText
    }
}

Replace-BulletText "Development of buildings internal load forecast models" "Development of the building’s internal load forecast models"
Replace-BulletText "Work is in progress" "Pecan Street data analysis done. The use of appliances is highly uncorrelated, hence not interesting for prediction."
Replace-BulletText "Development of Intelligent (RL) controllers for single building control" "Development of Intelligent (RL) controllers for single-building control"
Replace-BulletText "Work is in progress" "Dynamical system model with temporal logic constraints."

$tr = $tf.TextRange
$tr.InsertAfter("`rCan be tested on any dynamical system.`rOff-grid house model – RL with Temporal Logic constraints, including safety constraints.`rStart from the off-grid model and then proceed to grid-connected model flexibility provisioning.")

# Resize the textbox last, after all text edits (the box has spAutoFit,
# so resizing before editing the text would just get clobbered again).
$shp.Height = 312.62346456692916
